$d = $word.ActiveDocument

# --- Step 1: Re-type the first paragraph's text so Word collapses the
#     existing multi-run / proofErr-wrapped text into a single clean run
#     (this matches Word's normal behaviour when text is edited/retyped).
$firstParaText = "The iOS application is being targeted to iOS 7.0+. This is due to the fact that anything less than 7 does not support the core framework required for the Bluetooth functionality."
$find = $d.Content.Find
$find.Execute($firstParaText, $false, $false, $false, $false, $false, $true, 1, $false, $firstParaText, 2) | Out-Null

# --- Step 2: Add two blank paragraphs right before the paragraph that
#     carries the _GoBack bookmark (currently paragraph 2).
$bookmarkPara = $d.Paragraphs.Item(2)
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara = $d.Paragraphs.Item(3)
$bookmarkPara.Range.InsertParagraphBefore()

# --- Step 3: Type the new sentence at the very start of the bookmark
#     paragraph (now paragraph 4), ahead of the bookmark markers.
$bookmarkPara = $d.Paragraphs.Item(4)
$insertionPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$newSentence = "Designed for iPhones " + [char]0x2013 + " Isnt designed to be used for iPads as drivers don" + [char]0x2019 + "t generally use iPads when they are out driving " + [char]0x2013 + " RESEARCH THIS!!!"
$insertionPoint.InsertBefore($newSentence)

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
